# Update March 29, 2020
# Adds a new "29/03/2020" date column (M) to the DIY recovered-patients
# table, mirroring the existing date columns (C:L).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column L (header in row 2, data rows 3-81, and the blank filler rows
# 82-100) into the brand-new column M so every cell inherits the exact same
# formatting (borders, fonts, fills, number formats) already used by the
# other date columns.
$ws.Range("L2:L100").Copy($ws.Range("M2:M100"))

# The new column is for 29/03/2020.
$ws.Range("M2").Value = "29/03/2020"

# Data rows (3-81) get a numeric 0 value, same starting value as every other
# date column in the sheet.
$ws.Range("M3:M81").Value = 0

# Give column M the same width as the neighbouring date columns (K:L).
$ws.Range("M1").EntireColumn.ColumnWidth = $ws.Range("L1").EntireColumn.ColumnWidth
